$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for rows 2-51 per latest data pull
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.611.89"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.890.59"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.31"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4908"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2943"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06700"
$ws.Range("E9").Value = "  +1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.911.91"
$ws.Range("E10").Value = "  +1.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.15"
$ws.Range("E11").Value = "  +2.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07362"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.136"
$ws.Range("E13").Value = "  +2.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.17"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6681"
$ws.Range("E15").Value = "  +0.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.584.14"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007868"
$ws.Range("E17").Value = "  +1.26%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.41"
$ws.Range("E18").Value = "  +4.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.149.86"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.326"
$ws.Range("E21").Value = "  +12.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.004"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "189.51"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.200"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.510"
$ws.Range("E25").Value = "  +3.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.09"
$ws.Range("E26").Value = "  +2.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.931"
$ws.Range("E28").Value = "  +4.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.456"
$ws.Range("E29").Value = "  +4.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.391"
$ws.Range("E30").Value = "  +3.76%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09142"
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.065"
$ws.Range("E32").Value = "  +3.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05244"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7415"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("E35").Value = "  +2.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.719"
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01830"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.691"
$ws.Range("E38").Value = "  +1.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9206"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.063"
$ws.Range("E40").Value = "  +1.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "74.42"
$ws.Range("E41").Value = "  +31.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4421"
$ws.Range("E42").Value = "  +2.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.923"
$ws.Range("E43").Value = "  +5.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "106.09"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9940"
$ws.Range("E45").Value = "  -0.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1381"
$ws.Range("E46").Value = "  +4.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.553"
$ws.Range("E47").Value = "  +4.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05837"
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("E51").Value = "  +1.90%  "

# Rows 48-49: EnergySwap/Elrond swapped ranking order plus refreshed data
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.041"
$ws.Range("E48").Value = "  +5.10%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.21"
$ws.Range("E49").Value = "  +6.63%  "
